# "9th Stab - Cosmetic Changes"
# The watchlist gains two newer date columns (Jun_15, Jun_17) inserted in
# front of the existing Jun_13 / Jun_10 columns, and BidaskClub's tracker
# picks up a brand-new upgrade note for 6/16/2018, highlighted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh columns right before the old column B ("Jun_13").
# Everything that used to live in B/C slides over to D/E automatically,
# formulas and all.
$ws.Columns("B:C").Insert()

# New header row: newest dates go first (B=Jun_17, C=Jun_15), matching the
# existing newest-first ordering of D (Jun_13) / E (Jun_10).
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# The two new columns start out as "no rating change" (UN) placeholders for
# every analyst row, just like the other date columns.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# BidaskClub (row 22) actually has a new rating action on 6/16/2018 -
# record it in the newest column and highlight the cell so it stands out.
$bidaskCell = $ws.Cells.Item(22, 2)
$bidaskCell.Value = "6/16/2018,Upgrades,Sell -> Hold,"
$bidaskCell.Interior.ColorIndex = 42

# Keep the date columns at the same cosmetic width as before (8 chars).
$ws.Columns("C").ColumnWidth = 7.14
$ws.Columns("D").ColumnWidth = 7.14
$ws.Columns("E").ColumnWidth = 7.14
